# backtester/test/calculations.xlsx
#
# The backtester refactor split the "backtests" that used to live at
# columns I:J (AVG / SD / SHARPE labels + values) apart from the rest of
# the sheet by inserting a fresh, empty column at I. Everything that used
# to live in columns I/J shifts one column to the right (I -> J, J -> K),
# and Excel's formula engine automatically re-writes every relative
# reference that pointed into the old J column (e.g. SQRT(252)*J5/J6) so
# it keeps pointing at the same (now shifted) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("I").Insert()

# The author's selection/view also moved: the sheet no longer pins the
# top-left visible cell, and the active cell is now B4 instead of E2.
$ws.Range("B4").Select()
